$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1439.421
$ws.Range("I19").Value = 805.1429000000001
$ws.Range("J19").Value = 1809.4166
$ws.Range("K19").Value = 805.1429000000001
$ws.Range("L19").Value = 1809.4166
$ws.Range("M19").Value = -630.1429000000001
$ws.Range("N19").Value = -2159.4166

$ws.Range("H33").Value = 174.29411
$ws.Range("I33").Value = 160.26666
$ws.Range("J33").Value = 279.5
$ws.Range("K33").Value = 160.26666
$ws.Range("L33").Value = 279.5
$ws.Range("M33").Value = 68.73334
$ws.Range("N33").Value = -737.5

$ws.Range("H34").Value = 35400
$ws.Range("I34").Value = 13500
$ws.Range("K34").Value = 13500
$ws.Range("M34").Value = -13297

$ws.Range("H36").Value = 35400
$ws.Range("I36").Value = 13500
$ws.Range("K36").Value = 13500
$ws.Range("M36").Value = -12785

$ws.Range("H64").Value = 1894900
$ws.Range("J64").Value = 4794.727
$ws.Range("L64").Value = 4794.727
$ws.Range("N64").Value = -5290.727

$ws.Range("H67").Value = 1894900
$ws.Range("J67").Value = 4794.727
$ws.Range("L67").Value = 4794.727
$ws.Range("N67").Value = -6510.727

$ws.Range("H74").Value = 5000.6
$ws.Range("I74").Value = 5001
$ws.Range("K74").Value = 5001
$ws.Range("M74").Value = -4065

$ws.Range("H77").Value = 5000.6
$ws.Range("I77").Value = 5001
$ws.Range("K77").Value = 25005
$ws.Range("M77").Value = -20325

$ws.Range("H80").Value = 20833866
$ws.Range("I80").Value = 35714464
$ws.Range("K80").Value = 107143392
$ws.Range("M80").Value = -107142394

$ws.Range("H83").Value = 20833866
$ws.Range("I83").Value = 35714464
$ws.Range("K83").Value = 321430176
$ws.Range("M83").Value = -321425184

$ws.Range("H88").Value = 33344916
$ws.Range("I88").Value = 66670170
$ws.Range("J88").Value = 19666
$ws.Range("K88").Value = 66670170
$ws.Range("L88").Value = 19666
$ws.Range("M88").Value = -66669764
$ws.Range("N88").Value = -20478

$ws.Range("H91").Value = 33344916
$ws.Range("I91").Value = 66670170
$ws.Range("J91").Value = 19666
$ws.Range("K91").Value = 66670170
$ws.Range("L91").Value = 19666
$ws.Range("M91").Value = -66668766
$ws.Range("N91").Value = -22474

$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2504
$ws.Range("N97").ClearContents()

$ws.Range("H111").Value = 4603.8
$ws.Range("I111").Value = 4491.857
$ws.Range("K111").Value = 13475.571
$ws.Range("M111").Value = -10408.571

$ws.Range("H132").Value = 9036.17
$ws.Range("I132").Value = 3407.625
$ws.Range("K132").Value = 10222.875
$ws.Range("M132").Value = -7692.875

$ws.Range("H135").Value = 2784.7812
$ws.Range("J135").Value = 4132.933
$ws.Range("L135").Value = 37196.397
$ws.Range("N135").Value = -42266.397

$ws.Range("H138").Value = 8078.196
$ws.Range("I138").Value = 2999.2856
$ws.Range("J138").Value = 8989.795
$ws.Range("K138").Value = 8997.856800000001
$ws.Range("L138").Value = 26969.385
$ws.Range("M138").Value = -3857.856800000001
$ws.Range("N138").Value = -37249.385

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2270.3044
$ws.Range("I45").Value = 1195.2941
$ws.Range("K45").Value = 1195.2941
$ws.Range("M45").Value = -818.2941000000001

$ws.Range("H63").Value = 3223.25
$ws.Range("J63").Value = 2699
$ws.Range("L63").Value = 2699
$ws.Range("N63").Value = -4071

$ws.Range("H66").Value = 3223.25
$ws.Range("J66").Value = 2699
$ws.Range("L66").Value = 13495
$ws.Range("N66").Value = -20359

$ws.Range("H110").Value = 1278776.2
$ws.Range("I110").Value = 1856339.2
$ws.Range("K110").Value = 1856339.2
$ws.Range("M110").Value = -1854294.2

$ws.Range("H122").Value = 618722
$ws.Range("I122").Value = 1378624.5
$ws.Range("K122").Value = 4135873.5
$ws.Range("M122").Value = -4133423.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3506.5107
$ws.Range("I134").Value = 2239.8
$ws.Range("J134").Value = 7201.0835
$ws.Range("K134").Value = 6719.400000000001
$ws.Range("L134").Value = 21603.2505
$ws.Range("M134").Value = -4184.400000000001
$ws.Range("N134").Value = -26673.2505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 39462
$ws.Range("I62").Value = 2950
$ws.Range("J62").Value = 51632.668
$ws.Range("K62").Value = 2950
$ws.Range("L62").Value = 51632.668
$ws.Range("M62").Value = -2326
$ws.Range("N62").Value = -52880.668

$ws.Range("H65").Value = 39462
$ws.Range("I65").Value = 2950
$ws.Range("J65").Value = 51632.668
$ws.Range("K65").Value = 14750
$ws.Range("L65").Value = 258163.34
$ws.Range("M65").Value = -11630
$ws.Range("N65").Value = -264403.34

$ws.Range("H132").Value = 95255380
$ws.Range("I132").Value = 121215576
$ws.Range("K132").Value = 363646728
$ws.Range("M132").Value = -363644198

$ws.Range("H134").Value = 1599.6666
$ws.Range("I134").Value = 899.5
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 2698.5
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -163.5
$ws.Range("N134").Value = -14070

$ws.Range("H135").Value = 63166.582
$ws.Range("J135").Value = 63166.582
$ws.Range("L135").Value = 63166.582
$ws.Range("N135").Value = -73306.58199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 850.9524
$ws.Range("I122").Value = 693
$ws.Range("J122").Value = 1107.625
$ws.Range("K122").Value = 6237
$ws.Range("L122").Value = 9968.625
$ws.Range("M122").Value = -3787
$ws.Range("N122").Value = -14868.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 20839276
$ws.Range("I113").Value = 33338672
$ws.Range("K113").Value = 33338672
$ws.Range("M113").Value = -33336502

$ws.Range("H132").Value = 2995.425
$ws.Range("I132").Value = 3033.9062
$ws.Range("J132").Value = 2841.5
$ws.Range("K132").Value = 9101.7186
$ws.Range("L132").Value = 8524.5
$ws.Range("M132").Value = -6571.7186
$ws.Range("N132").Value = -13584.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8950.571
$ws.Range("I16").Value = 7468.4707
$ws.Range("J16").Value = 15249.5
$ws.Range("K16").Value = 7468.4707
$ws.Range("L16").Value = 15249.5
$ws.Range("M16").Value = -7298.4707
$ws.Range("N16").Value = -15589.5

$ws.Range("H22").Value = 1842
$ws.Range("I22").Value = 1414.3334
$ws.Range("J22").Value = 2611.8
$ws.Range("K22").Value = 1414.3334
$ws.Range("L22").Value = 2611.8
$ws.Range("M22").Value = -1119.3334
$ws.Range("N22").Value = -3201.8

$ws.Range("H27").Value = 1842
$ws.Range("I27").Value = 1414.3334
$ws.Range("J27").Value = 2611.8
$ws.Range("K27").Value = 1414.3334
$ws.Range("L27").Value = 2611.8
$ws.Range("M27").Value = -1307.3334
$ws.Range("N27").Value = -2825.8

$ws.Range("H50").Value = 13443.333
$ws.Range("I50").Value = 12998.571
$ws.Range("J50").Value = 15000
$ws.Range("K50").Value = 12998.571
$ws.Range("L50").Value = 15000
$ws.Range("M50").Value = -12361.571
$ws.Range("N50").Value = -16274

$ws.Range("H55").Value = 247.1875
$ws.Range("I55").Value = 170.6
$ws.Range("K55").Value = 170.6
$ws.Range("M55").Value = 2.400000000000006

$ws.Range("H100").Value = 6999.75
$ws.Range("I100").Value = 5999.6665
$ws.Range("K100").Value = 5999.6665
$ws.Range("M100").Value = -5458.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1171.4445
$ws.Range("I107").Value = 778.5
$ws.Range("J107").Value = 2546.75
$ws.Range("K107").Value = 2335.5
$ws.Range("L107").Value = 7640.25
$ws.Range("M107").Value = -415.5
$ws.Range("N107").Value = -11480.25

$ws.Range("H132").Value = 26325160
$ws.Range("I132").Value = 9973.177
$ws.Range("J132").Value = 250004260
$ws.Range("K132").Value = 29919.531
$ws.Range("L132").Value = 750012780
$ws.Range("M132").Value = -27389.531
$ws.Range("N132").Value = -750017840

$ws.Range("H137").Value = 100000
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200
